# Forecast Comparison sheet: shift each week's Week_Start_Date forward by
# one week and update MyForecast (column D) with the new values.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$newDates = @(
    "2025-01-12", "2025-01-19", "2025-01-26", "2025-02-02",
    "2025-02-09", "2025-02-16", "2025-02-23", "2025-03-02",
    "2025-03-09", "2025-03-16", "2025-03-23", "2025-03-30",
    "2025-04-06", "2025-04-13", "2025-04-20", "2025-04-27"
)

$newForecast = @(57, 58, 60, 62, 64, 66, 65, 66, 65, 63, 59, 52, 49, 47, 46, 46)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $i + 2

    # Column B holds the date as plain text in the source file (inline
    # string), not a real Excel date. Force text entry via NumberFormat="@"
    # so the date-looking string isn't auto-converted into a date serial,
    # then clear the format again so no residual style is left behind.
    $cellB = $ws1.Cells.Item($row, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $newDates[$i]
    $cellB.ClearFormats()

    # Column D (MyForecast) is a genuine number in the source file.
    $ws1.Cells.Item($row, 4).Value = $newForecast[$i]
}

# Summary sheet: update the recomputed aggregate metrics. All "Value"
# column cells on this sheet are stored as text in the source file (even
# the numeric-looking ones), so every write here uses the same
# NumberFormat="@" / ClearFormats() trick to keep them as plain text
# without leaving a residual style behind.
$ws2 = $wb.Worksheets.Item("Summary")

function Set-TextValue($sheet, $row, $col, $value) {
    $cell = $sheet.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws2 2 2 "2023-01-22 to 2025-01-05"   # Historical Range
Set-TextValue $ws2 4 2 "121"                         # Max Sales
Set-TextValue $ws2 6 2 "42"                          # Median Sales
Set-TextValue $ws2 8 2 "4309 units"                  # Total Historical Sales
Set-TextValue $ws2 9 2 "925"                         # Total Forecast (16 Weeks)
Set-TextValue $ws2 10 2 "498"                        # Total Forecast (8 Weeks)
Set-TextValue $ws2 11 2 "237"                        # Total Forecast (4 Weeks)
Set-TextValue $ws2 12 2 "66"                         # Max Forecast
Set-TextValue $ws2 13 2 "2025-02-16"                 # Max Forecast Week
Set-TextValue $ws2 14 2 "46"                         # Min Forecast
Set-TextValue $ws2 15 2 "2025-04-20"                 # Min Forecast Week
